# Update DiffExcel to V1.2
# Adds a fourth scenario (Scenario4 / Sheet4) results column (column F) to
# the "DiffWorksheet" sheet - mirroring the existing Scenario1/2/3 columns
# (C/D/E) - and makes "DiffWorksheet" the active sheet/tab instead of
# "DiffExcel".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("DiffExcel")
$ws2 = $wb.Worksheets.Item("DiffWorksheet")

# --- Give the new column F the same look (style/width) as column E -----
$ws2.Range("E4:E13").Copy()
$ws2.Range("F4:F13").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Columns.Item(6).ColumnWidth = $ws2.Columns.Item(5).ColumnWidth

# --- Populate the new "Scenario4" column --------------------------------
$ws2.Range("F4").Value  = "Scenario4"
$ws2.Range("F5").Value  = "Sheet4"
$ws2.Range("F6").Value  = "Object[]"
$ws2.Range("F7").Value  = 3
$ws2.Range("F8").Value  = "0"
$ws2.Range("F9").Value  = "A1   OldA1"
$ws2.Range("F10").Value = "1"
$ws2.Range("F11").Value = "B1   OldB1"
$ws2.Range("F12").Value = "2"
$ws2.Range("F13").Value = "A2   OldA2"

# --- Make "DiffWorksheet" the selected/active sheet ---------------------
$ws2.Select()
$ws2.Range("F14").Select()
$wb.Windows.Item(1).ActiveSheet = $ws2
